# Update GAM, NGAM Figure growth
# - Insert a new "Strain" column (C) into the RateCompare sheet, shifting the
#   existing BMLaw..Yield-g/g columns one place to the right (D..N).
# - Populate the new column with the strain identifier for every data row
#   (numeric "521" for most strains, text "MB215" for the two Wierckx rows).
# - The new column is formatted as Text, matching the rest of the sheet's
#   informal "as typed" styling.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RateCompare")

# Insert a new blank column before the existing column C ("BMLaw" and
# everything to its right shifts one column to the right, C -> D, ... M -> N).
$ws.Columns("C").Insert()

# Header
$ws.Range("C1").Value = "Strain"

# Strain id per row. Numbers are written as real numbers first so the stored
# cell type stays numeric; the Text number format is applied afterwards so it
# only affects display/formatting, not the underlying stored type (matches
# rows where the value is numeric, e.g. 521, and rows where it is the literal
# text "MB215").
$strainValues = @{
    2  = 521
    3  = 521
    4  = 521
    5  = 521
    6  = 521
    7  = "MB215"
    8  = "MB215"
    9  = 521
    10 = 521
}

foreach ($row in $strainValues.Keys) {
    $ws.Range("C$row").Value = $strainValues[$row]
}

# Apply the "Text" number format to the whole new column, like the rest of
# the sheet's string-typed columns.
$ws.Columns("C").NumberFormat = "@"

# Restore the user's final selection on this sheet.
$ws.Activate()
$ws.Range("D22").Select()
